# Atualiza o roteiro de entregas: corrige a lógica para não gerar a mesma
# rota 2 vezes, reordenando/atualizando os pedidos de 13-15/11/2024 e
# adicionando os pedidos que faltavam (linhas 15-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linha modelo (formatação das linhas de dados já existente no arquivo).
$templateRow = 2

# Garante que as novas linhas (15-20) herdem a mesma formatação (bordas,
# alinhamento, fonte) das linhas de dados já existentes, copiando apenas o
# formato antes de escrever qualquer valor.
$ws.Range("A$templateRow`:E$templateRow").Copy()
$ws.Range("A15:E20").PasteSpecial(-4122)

# Dados finais (linha, Nº Pedido, Data de entrega, Período, Nome do cliente, Bairro)
$data = @(
    @(2, '11069', '13/11/2024', 'manhã', 'Desconhecido', 'Barreiros'),
    @(3, '11068', '13/11/2024', 'manhã', 'Desconhecido', 'Barreiros'),
    @(4, '11074', '13/11/2024', 'manhã', 'Desconhecido', 'Bela Vista'),
    @(5, '11107', '13/11/2024', 'manhã', 'Desconhecido', 'Campinas'),
    @(6, '11084', '13/11/2024', 'tarde', 'Desconhecido', 'Brejarú'),
    @(7, '11083', '13/11/2024', 'tarde', 'Desconhecido', 'Brejarú'),
    @(8, '11076', '13/11/2024', 'tarde', 'Desconhecido', 'Bela Vista'),
    @(9, '11061', '13/11/2024', 'tarde', 'Desconhecido', 'Aririú'),
    @(10, '11078', '13/11/2024', 'tarde', 'Desconhecido', 'Bela Vista'),
    @(11, '11080', '14/11/2024', 'manhã', 'Desconhecido', 'Bela Vista'),
    @(12, '11077', '14/11/2024', 'manhã', 'Desconhecido', 'Bela Vista'),
    @(13, '11057', '14/11/2024', 'manhã', 'Desconhecido', 'Aririú'),
    @(14, '11114', '14/11/2024', 'tarde', 'Desconhecido', 'Campinas'),
    @(15, '11110', '14/11/2024', 'tarde', 'Desconhecido', 'Campinas'),
    @(16, '11111', '14/11/2024', 'tarde', 'Desconhecido', 'Campinas'),
    @(17, '11115', '14/11/2024', 'tarde', 'Desconhecido', 'Campinas'),
    @(18, '11050', '15/11/2024', 'manhã', 'Desconhecido', 'Areias de Cima (Guaporanga)'),
    @(19, '11088', '15/11/2024', 'manhã', 'Desconhecido', 'Cachoeiras (Guaporanga)'),
    @(20, '11105', '15/11/2024', 'tarde', 'Desconhecido', 'Campeche')
)

foreach ($row in $data) {
    $r = $row[0]
    $pedido = $row[1]
    $dataEntrega = $row[2]
    $periodo = $row[3]
    $cliente = $row[4]
    $bairro = $row[5]

    # Nº Pedido é puramente numérico: usa apóstrofo para forçar o
    # armazenamento como texto (preserva zeros à esquerda / tipo string),
    # igual ao arquivo original (inlineStr).
    $ws.Cells.Item($r, 1).Value2 = "'" + $pedido
    $ws.Cells.Item($r, 2).Value2 = $dataEntrega
    $ws.Cells.Item($r, 3).Value2 = $periodo
    $ws.Cells.Item($r, 4).Value2 = $cliente
    $ws.Cells.Item($r, 5).Value2 = $bairro
}
